$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1306.8572
$ws.Range("J19").Value = 1406.3334
$ws.Range("L19").Value = 1406.3334
$ws.Range("N19").Value = -1756.3334
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 100
$ws.Range("K20").Value = 100
$ws.Range("M20").Value = 130
$ws.Range("H31").Value = 2126.5
$ws.Range("I31").Value = 117.71429
$ws.Range("K31").Value = 353.14287
$ws.Range("M31").Value = -123.14287
$ws.Range("H35").Value = 100
$ws.Range("I35").Value = 100
$ws.Range("K35").Value = 100
$ws.Range("M35").Value = 279
$ws.Range("H62").Value = 62501936
$ws.Range("I62").Value = 83334910
$ws.Range("K62").Value = 83334910
$ws.Range("M62").Value = -83334286
$ws.Range("H65").Value = 62501936
$ws.Range("I65").Value = 83334910
$ws.Range("K65").Value = 416674550
$ws.Range("M65").Value = -416671430
$ws.Range("H70").Value = 4833.6875
$ws.Range("I70").Value = 1539
$ws.Range("J70").Value = 16600.428
$ws.Range("K70").Value = 4617
$ws.Range("L70").Value = 49801.284
$ws.Range("M70").Value = -4347
$ws.Range("N70").Value = -50341.284
$ws.Range("H73").Value = 4833.6875
$ws.Range("I73").Value = 1539
$ws.Range("J73").Value = 16600.428
$ws.Range("K73").Value = 4617
$ws.Range("L73").Value = 49801.284
$ws.Range("M73").Value = -3681
$ws.Range("N73").Value = -51673.284
$ws.Range("H76").Value = 71432620
$ws.Range("I76").Value = 200003360
$ws.Range("J76").Value = 4431.5557
$ws.Range("K76").Value = 200003360
$ws.Range("L76").Value = 4431.5557
$ws.Range("M76").Value = -200003045
$ws.Range("N76").Value = -5061.5557
$ws.Range("H79").Value = 71432620
$ws.Range("I79").Value = 200003360
$ws.Range("J79").Value = 4431.5557
$ws.Range("K79").Value = 200003360
$ws.Range("L79").Value = 4431.5557
$ws.Range("M79").Value = -200002268
$ws.Range("N79").Value = -6615.5557
$ws.Range("H98").Value = 1652.3636
$ws.Range("I98").Value = 1568.7142
$ws.Range("J98").Value = 1798.75
$ws.Range("K98").Value = 1568.7142
$ws.Range("L98").Value = 1798.75
$ws.Range("M98").Value = -70.71419999999989
$ws.Range("N98").Value = -4794.75
$ws.Range("H122").Value = 1652.3636
$ws.Range("I122").Value = 1568.7142
$ws.Range("J122").Value = 1798.75
$ws.Range("K122").Value = 4706.142599999999
$ws.Range("L122").Value = 5396.25
$ws.Range("M122").Value = -2256.142599999999
$ws.Range("N122").Value = -10296.25
$ws.Range("H131").Value = 3972.125
$ws.Range("I131").Value = 2411
$ws.Range("K131").Value = 7233
$ws.Range("M131").Value = -2193

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 209.57143
$ws.Range("I5").Value = 146.75
$ws.Range("K5").Value = 146.75
$ws.Range("M5").Value = -34.75
$ws.Range("H132").Value = 1656.25
$ws.Range("I132").Value = 1491.3715
$ws.Range("J132").Value = 2810.4
$ws.Range("K132").Value = 4474.1145
$ws.Range("L132").Value = 8431.200000000001
$ws.Range("M132").Value = -1944.1145
$ws.Range("N132").Value = -13491.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 209.57143
$ws.Range("I4").Value = 146.75
$ws.Range("K4").Value = 146.75
$ws.Range("M4").Value = -31.75
$ws.Range("H134").Value = 6040.018
$ws.Range("I134").Value = 4772.3257
$ws.Range("J134").Value = 10233.154
$ws.Range("K134").Value = 14316.9771
$ws.Range("L134").Value = 30699.462
$ws.Range("M134").Value = -11781.9771
$ws.Range("N134").Value = -35769.462
$ws.Range("H139").Value = 197166.5
$ws.Range("J139").Value = 197166.5
$ws.Range("L139").Value = 197166.5
$ws.Range("N139").Value = -207446.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1901.1052
$ws.Range("I99").Value = 1765.6364
$ws.Range("K99").Value = 1765.6364
$ws.Range("M99").Value = -267.6364000000001
$ws.Range("H126").Value = 1901.1052
$ws.Range("I126").Value = 1765.6364
$ws.Range("K126").Value = 5296.9092
$ws.Range("M126").Value = -2826.9092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 551.3125
$ws.Range("I8").Value = 551.3125
$ws.Range("K8").Value = 1653.9375
$ws.Range("M8").Value = -1514.9375
$ws.Range("H39").Value = 6225.75
$ws.Range("J39").Value = 8234.333000000001
$ws.Range("L39").Value = 24702.999
$ws.Range("N39").Value = -25290.999
$ws.Range("H107").Value = 658.5333000000001
$ws.Range("I107").Value = 261.6
$ws.Range("K107").Value = 784.8000000000001
$ws.Range("M107").Value = 1135.2
$ws.Range("H121").Value = 81815.266
$ws.Range("I121").Value = 500
$ws.Range("J121").Value = 87623.5
$ws.Range("K121").Value = 1500
$ws.Range("L121").Value = 262870.5
$ws.Range("M121").Value = -190
$ws.Range("N121").Value = -265490.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9484.565000000001
$ws.Range("I70").Value = 10255.5
$ws.Range("J70").Value = 7722.4287
$ws.Range("K70").Value = 10255.5
$ws.Range("L70").Value = 7722.4287
$ws.Range("M70").Value = -9985.5
$ws.Range("N70").Value = -8262.4287
$ws.Range("H73").Value = 9484.565000000001
$ws.Range("I73").Value = 10255.5
$ws.Range("J73").Value = 7722.4287
$ws.Range("K73").Value = 10255.5
$ws.Range("L73").Value = 7722.4287
$ws.Range("M73").Value = -9319.5
$ws.Range("N73").Value = -9594.4287
$ws.Range("H141").Value = 128236.875
$ws.Range("J141").Value = 132270.72
$ws.Range("L141").Value = 132270.72
$ws.Range("N141").Value = -142630.72

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 458038.34
$ws.Range("I20").Value = 77825.914
$ws.Range("K20").Value = 77825.914
$ws.Range("M20").Value = -77599.914
$ws.Range("H122").Value = 6307.737
$ws.Range("I122").Value = 4699.25
$ws.Range("J122").Value = 6736.6665
$ws.Range("K122").Value = 14097.75
$ws.Range("L122").Value = 20209.9995
$ws.Range("M122").Value = -11647.75
$ws.Range("N122").Value = -25109.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3373.3225
$ws.Range("I122").Value = 2646.1155
$ws.Range("J122").Value = 7154.8
$ws.Range("K122").Value = 7938.3465
$ws.Range("L122").Value = 21464.4
$ws.Range("M122").Value = -5488.3465
$ws.Range("N122").Value = -26364.4
$ws.Range("H126").Value = 2092
$ws.Range("I126").Value = 1997.2106
$ws.Range("K126").Value = 5991.6318
$ws.Range("M126").Value = -3521.6318
